$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.796.58"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.667.03"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.34"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.36"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.662"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +7.12%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.401"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.84"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.83"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000195"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.144.71"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.643.28"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.675.41"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.60"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.78"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.11"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.44"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.10%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.68"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.67%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.55"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.00%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "564.34"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.47%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.02"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.11"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.22%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.58"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.86%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.422"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.53"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.48"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "160.85"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.53%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.78%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0605"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.74"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.42%  "

$ws.Range("B46").Value = "Stellar"

$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.104"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.05%  "

$ws.Range("B47").Value = "Mantle"

$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.639"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.66%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.81"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0245"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.805"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.76%  "
